$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Before insert: B2=$($ws.Cells.Item(2,2).Text) C2=$($ws.Cells.Item(2,3).Text)"
$ws.Columns.Item(2).Insert()
Write-Host "After insert: B2=$($ws.Cells.Item(2,2).Text) C2=$($ws.Cells.Item(2,3).Text) D2=$($ws.Cells.Item(2,4).Text)"
